$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the "mean annual precipitation (mm)" values in column AI (col 35):
#    the values were recomputed and are 1.2x the previously stored numbers,
#    rows 2 through 166 (row 1 is the header).
for ($r = 2; $r -le 166; $r++) {
    $cell = $ws.Cells.Item($r, 35)
    $old = $cell.Value2
    $cell.Value = $old * 1.2
}

# 2. The AI column values should no longer carry the two-decimal number
#    format (style index 5 in the original file) - they revert to the
#    default/general style, matching an already-unstyled cell (e.g. A2).
$plainStyle = $ws.Range("A2").Style
$ws.Range("AI2:AI166").Style = $plainStyle

# 3. Column AI width shrinks slightly (30.75 -> 30.25 characters).
$ws.Columns("AI").ColumnWidth = 29.42

# 4. Update the active selection/cursor position to AJ6 (previously the
#    whole of column A was selected).
$ws.Range("AJ6").Select()
